# Auto-generated edit script: update crypto price/volume table to latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.792.91"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "1.781.51"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5130"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3781"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07773"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.085"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.194"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.60%  "
$ws.Range("D15").Value = "1.770.08"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.157"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001070"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.906"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").Value = "27.839.17"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.233"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("D28").Value = "1.983.04"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.347"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1077"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.476"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07042"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.87%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.664"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2117"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.90%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.015"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.76%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6077"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.148"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.321"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5954"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.894"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06703"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.06%  "
